$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header text updates (A1 / K1).
# K1 is set first, then A1, so that the shared-string table ends up
# compacted in the same order as the target workbook
# (index 13 = "AiHa Sort vs Counting Sort",
#  index 14 = "Num Integers/SortAlgorithm (Clicks)").
# ------------------------------------------------------------------
$ws.Range("K1").Value = "AiHa Sort vs Counting Sort"
$ws.Range("A1").Value = "Num Integers/SortAlgorithm (Clicks)"

# ------------------------------------------------------------------
# Column width updates.
# ColumnWidth is stored internally as (value + 5/6) rounded to the
# nearest 1/6, so we back-solve for the ColumnWidth that reproduces
# the target stored widths as closely as possible.
# Column A: 32.33203125 -> 41.33203125
# Column K: 25.83203125 -> 31.1640625
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 40.498697916666664
$ws.Columns.Item(11).ColumnWidth = 30.330729166666668

# ------------------------------------------------------------------
# Selection change: active cell moves from K6 to L12.
# ------------------------------------------------------------------
$ws.Range("L12").Select()
